$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 15: attachment for the PCR file flanks_short.xlsx, linked back to its
# SynBioHub record via a hyperlink on the display_id cell.
$ws.Hyperlinks.Add($ws.Range("A15"), "https://synbiohub.org/user/jhay/Johnny190421/sl0199_flatten/1.0.0", "", "", "https://synbiohub.org/user/jhay/Johnny190421/sl0199_flatten/1.0.0")
$ws.Range("A15").Value = "sl0199_flatten"

$ws.Range("B15").Value = "flanks_short.xlsx"
$ws.Range("D15").Formula = "=_xlfn.CONCAT(B15, "" is a PCR file"")"

$ws.Range("A15").Select()
